$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF") -- match the existing header
# formatting (bold font, thin border, center/top alignment) by copying the
# format from the neighboring header cell H1 rather than re-declaring it.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows 2-32: new values for columns I and J.
$data = @(
    @(8, 8),
    @(6, 7),
    @(6, 6),
    @(7, 7),
    @(10, 11),
    @(5, 5),
    @(6, 7),
    @(5, 5),
    @(10, 10),
    @(5, 6),
    @(6, 7),
    @(5, 6),
    @(11, 11),
    @(6, 7),
    @(3, 4),
    @(7, 8),
    @(5, 6),
    @(7, 7),
    @(7, 7),
    @(4, 5),
    @(5, 7),
    @(5, 6),
    @(8, 8),
    @(6, 7),
    @(7, 7),
    @(7, 8),
    @(4, 4),
    @(8, 8),
    @(6, 6),
    @(9, 9),
    @(7, 7)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
